# Applies the "Added year to archetypes" edit:
#  - ARCHITECTURE and HVAC sheets: fill the previously-empty year_start (B)
#    and year_end (C) columns (rows 2-19) with 1950 / 2030 respectively.
#    Both columns are formatted as Text, so the values are stored as text.
#  - INDOOR_COMFORT sheet: replace the placeholder "None" text used for the
#    "no heating system"/"no cooling system" set points with real numeric
#    values (10 for the low/no-heating set points in C & E, and 50/10 for
#    rows 13 & 15's B/C/D/E, which previously were entirely "None").

$wb = $excel.ActiveWorkbook

$archSheet = $wb.Worksheets.Item("ARCHITECTURE")
$hvacSheet = $wb.Worksheets.Item("HVAC")
$comfortSheet = $wb.Worksheets.Item("INDOOR_COMFORT")

# --- ARCHITECTURE: year_start / year_end ---------------------------------
for ($row = 2; $row -le 19; $row++) {
    $archSheet.Range("B$row").Value = "1950"
    $archSheet.Range("C$row").Value = "2030"
}
$archSheet.Range("B19:C19").Select()

# --- HVAC: year_start / year_end ------------------------------------------
for ($row = 2; $row -le 19; $row++) {
    $hvacSheet.Range("B$row").Value = "1950"
    $hvacSheet.Range("C$row").Value = "2030"
}
$hvacSheet.Range("B19:C19").Select()

# --- INDOOR_COMFORT: replace "None" placeholders with real numbers -------
for ($row = 2; $row -le 19; $row++) {
    $comfortSheet.Range("C$row").Value = 10
    $comfortSheet.Range("E$row").Value = 10
}

# Rows 13 (SWIMMING) and 15 (PARKING) had every set point as "None";
# give them real heating/cooling set points too.
foreach ($row in 13, 15) {
    $comfortSheet.Range("B$row").Value = 50
    $comfortSheet.Range("D$row").Value = 50
}

$comfortSheet.Range("H36").Select()

# Keep ARCHITECTURE as the active/visible tab, as in the original workbook.
$archSheet.Activate()
